$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "2022" column (N) is being added, mirroring the existing "2021"
# column (M): same per-row formatting (number format / font / borders),
# with new data values for each indicator row.

# Clone column M's formatting (rows 4-13, the header year row through the
# last data row) into column N so every new cell picks up the right
# number format, bold/italic font, and borders for its row.
$ws.Range("M4:M13").Copy() | Out-Null
$ws.Range("N4:N13").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# New "2022" header value.
$ws.Range("N4").Value = 2022

# New data for the 2022 column (one value per indicator row).
$ws.Range("N5").Value = 4.3
$ws.Range("N6").Value = 5.0999999999999996
$ws.Range("N7").Value = 3.1
$ws.Range("N8").Value = 2.9
$ws.Range("N9").Value = 3.4
$ws.Range("N10").Value = 2.2999999999999998
$ws.Range("N11").Value = 92.8
$ws.Range("N12").Value = 91.6
$ws.Range("N13").Value = 94.6

# The saved worksheet view now has N15 as the active selection.
$ws.Range("N15").Select() | Out-Null
